# API_driven.xlsx update — "update beforesuite keep extent report"
#
# Adds three new data rows (3,4,5) to Sheet2, mirroring the existing
# row 2 pattern (id / titile / dueDate / completed), extends the
# "completed" list-data-validation down to the new rows, and moves the
# active tab/selection from Sheet3 back to Sheet2 (Sheet3's own
# selection is left untouched).

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- New rows on Sheet2 --------------------------------------------------
# Column A uses the same "quote-prefixed text" look as A2 (e.g. '0), so a
# leading apostrophe reproduces that text-with-quote-prefix storage.
# Columns B/C reuse the textual values already used in row 2 (shared
# strings "string" / the ISO timestamp), forced to text format so they
# match row 2's cell style instead of being stored as numbers/dates.
# Column D is a plain boolean, same as D2.

$ws2.Range("A3").Value = "'1"
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "string"
$ws2.Range("C3").NumberFormat = "@"
$ws2.Range("C3").Value = "2023-02-10T11:01:06.626Z"
$ws2.Range("D3").Value = $true

$ws2.Range("A4").Value = "'2"
$ws2.Range("B4").NumberFormat = "@"
$ws2.Range("B4").Value = "string"
$ws2.Range("C4").NumberFormat = "@"
$ws2.Range("C4").Value = "2023-02-10T11:01:06.626Z"
$ws2.Range("D4").Value = $true

$ws2.Range("A5").Value = "'3"
$ws2.Range("B5").NumberFormat = "@"
$ws2.Range("B5").Value = "string"
$ws2.Range("C5").NumberFormat = "@"
$ws2.Range("C5").Value = "2023-02-10T11:01:06.626Z"
$ws2.Range("D5").Value = $false

# --- Extend the "completed" list validation from D2 to D2:D5 -------------
$ws2.Range("D2:D5").Validation.Delete()
$ws2.Range("D2:D5").Validation.Add(3, 1, 1, '"true,false"')
$ws2.Range("D2:D5").Validation.InCellDropdown = $true

# --- Selection / active sheet --------------------------------------------
# Originally Sheet3 was the active tab with its own K12 selection; the
# edit moves the active tab back to Sheet2 and leaves a new selection on
# it, while Sheet3 keeps its selection but is no longer the active tab.
$ws3.Range("K12").Select()
$ws2.Activate()
$ws2.Range("C16").Select()
